# Auto-generated Excel COM-interop edit script
# Applies updated currentAveragePrice / Leve profit figures across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17 (hunk 0)
$ws.Range("H17").Value = 2434.5715
$ws.Range("J17").Value = 2517.3704
$ws.Range("L17").Value = 7552.111199999999
$ws.Range("N17").Value = -7888.111199999999
# Row 53 (hunk 1)
$ws.Range("H53").Value = 173.5
$ws.Range("I53").Value = 59.625
$ws.Range("J53").Value = 264.6
$ws.Range("K53").Value = 59.625
$ws.Range("L53").Value = 264.6
$ws.Range("M53").Value = 577.375
$ws.Range("N53").Value = -1538.6
# Row 118 (hunk 2)
$ws.Range("H118").Value = 766.0217
$ws.Range("I118").Value = 396.29413
$ws.Range("J118").Value = 982.7586
$ws.Range("K118").Value = 1188.88239
$ws.Range("L118").Value = 2948.2758
$ws.Range("M118").Value = 468.11761
$ws.Range("N118").Value = -6262.275799999999
# Row 132 (hunk 3)
$ws.Range("H132").Value = 1763.8136
$ws.Range("I132").Value = 1281.7174
$ws.Range("J132").Value = 3469.6924
$ws.Range("K132").Value = 3845.1522
$ws.Range("L132").Value = 10409.0772
$ws.Range("M132").Value = -1315.1522
$ws.Range("N132").Value = -15469.0772
# Row 135 (hunk 4)
$ws.Range("H135").Value = 360.40625
$ws.Range("I135").Value = 336.5862
$ws.Range("J135").Value = 590.6667
$ws.Range("K135").Value = 3029.2758
$ws.Range("L135").Value = 5316.0003
$ws.Range("M135").Value = -494.2758000000003
$ws.Range("N135").Value = -10386.0003
# Row 137 (hunk 5)
$ws.Range("H137").Value = 3243.5264
$ws.Range("I137").Value = 1080.2222
$ws.Range("K137").Value = 3240.6666
$ws.Range("M137").Value = -690.6665999999996
# Row 138 (hunk 6)
$ws.Range("H138").Value = 2718.6956
$ws.Range("I138").Value = 1351.6154
$ws.Range("J138").Value = 4495.9
$ws.Range("K138").Value = 4054.8462
$ws.Range("L138").Value = 13487.7
$ws.Range("M138").Value = 1085.1538
$ws.Range("N138").Value = -23767.7

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (hunk 7)
$ws.Range("H32").Value = 4606
$ws.Range("I32").Value = 4952.6978
$ws.Range("J32").Value = 2121.3333
$ws.Range("K32").Value = 4952.6978
$ws.Range("L32").Value = 2121.3333
$ws.Range("M32").Value = -4665.6978
$ws.Range("N32").Value = -2695.3333
# Row 34 (hunk 8)
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
# Row 74 (hunk 9)
$ws.Range("H74").Value = 3391.9556
$ws.Range("I74").Value = 1034.6487
$ws.Range("K74").Value = 1034.6487
$ws.Range("M74").Value = -160.6487
# Row 77 (hunk 10)
$ws.Range("H77").Value = 3391.9556
$ws.Range("I77").Value = 1034.6487
$ws.Range("K77").Value = 5173.2435
$ws.Range("M77").Value = -805.2434999999996
# Row 102 (hunk 11)
$ws.Range("H102").Value = 1905
$ws.Range("I102").Value = 1810
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1810
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -188
$ws.Range("N102").Value = -5244
# Row 132 (hunk 12)
$ws.Range("H132").Value = 6235.7617
$ws.Range("I132").Value = 4334.2583
$ws.Range("J132").Value = 11594.546
$ws.Range("K132").Value = 13002.7749
$ws.Range("L132").Value = 34783.638
$ws.Range("M132").Value = -10472.7749
$ws.Range("N132").Value = -39843.638

$ws = $wb.Worksheets.Item("BSM")
# Row 38 (hunk 13)
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 35 (hunk 14)
$ws.Range("H35").Value = 643.75
$ws.Range("I35").Value = 643.75
$ws.Range("K35").Value = 643.75
$ws.Range("M35").Value = -349.75
# Row 99 (hunk 15)
$ws.Range("H99").Value = 2829.3333
$ws.Range("I99").Value = 1712
$ws.Range("J99").Value = 3388
$ws.Range("K99").Value = 1712
$ws.Range("L99").Value = 3388
$ws.Range("M99").Value = -214
$ws.Range("N99").Value = -6384
# Row 126 (hunk 16)
$ws.Range("H126").Value = 2829.3333
$ws.Range("I126").Value = 1712
$ws.Range("J126").Value = 3388
$ws.Range("K126").Value = 5136
$ws.Range("L126").Value = 10164
$ws.Range("M126").Value = -2666
$ws.Range("N126").Value = -15104
# Row 134 (hunk 17)
$ws.Range("H134").Value = 867.8095
$ws.Range("I134").Value = 692
$ws.Range("K134").Value = 2076
$ws.Range("M134").Value = 459

$ws = $wb.Worksheets.Item("CUL")
# Row 36 (hunk 18)
$ws.Range("H36").Value = 55556410
$ws.Range("I36").Value = 1020.4
$ws.Range("K36").Value = 3061.2
$ws.Range("M36").Value = -2892.2
# Row 134 (hunk 19)
$ws.Range("H134").Value = 24279.66
$ws.Range("I134").Value = 59357
$ws.Range("J134").Value = 2780.6453
$ws.Range("K134").Value = 178071
$ws.Range("L134").Value = 8341.9359
$ws.Range("M134").Value = -173001
$ws.Range("N134").Value = -18481.9359

$ws = $wb.Worksheets.Item("GSM")
# Row 21 (hunk 20)
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
# Row 30 (hunk 21)
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
# Row 126 (hunk 22)
$ws.Range("H126").Value = 1220.2727
$ws.Range("I126").Value = 1073.2858
$ws.Range("J126").Value = 1477.5
$ws.Range("K126").Value = 3219.8574
$ws.Range("L126").Value = 4432.5
$ws.Range("M126").Value = -749.8574000000003
$ws.Range("N126").Value = -9372.5

$ws = $wb.Worksheets.Item("LTW")
# Row 32 (hunk 23)
$ws.Range("H32").Value = 782.4
$ws.Range("I32").Value = 782.4
$ws.Range("K32").Value = 782.4
$ws.Range("M32").Value = -465.4
# Row 46 (hunk 24)
$ws.Range("H46").Value = 1783.8823
$ws.Range("I46").Value = 2284
$ws.Range("J46").Value = 867
$ws.Range("K46").Value = 2284
$ws.Range("L46").Value = 867
$ws.Range("M46").Value = -2096
$ws.Range("N46").Value = -1243
# Row 132 (hunk 25)
$ws.Range("H132").Value = 6576.4614
$ws.Range("I132").Value = 8399.177
$ws.Range("J132").Value = 3133.5557
$ws.Range("K132").Value = 25197.531
$ws.Range("L132").Value = 9400.667099999999
$ws.Range("M132").Value = -22667.531
$ws.Range("N132").Value = -14460.6671

$ws = $wb.Worksheets.Item("WVR")
# Row 45 (hunk 26)
$ws.Range("H45").Value = 7811.3335
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 7811.3335
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 7811.3335
$ws.Range("N45").Value = -8793.333500000001
$ws.Range("M45").ClearContents()
# Row 132 (hunk 27)
$ws.Range("H132").Value = 8575.294
$ws.Range("I132").Value = 15139.625
$ws.Range("J132").Value = 2740.3333
$ws.Range("K132").Value = 45418.875
$ws.Range("L132").Value = 8220.999899999999
$ws.Range("M132").Value = -42888.875
$ws.Range("N132").Value = -13280.9999
